$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("B1").Value = 20
$ws.Range("F1").Value = "A"
$ws.Range("G1").Value = "A"
$ws.Range("H1").Value = "A"
$ws.Range("I1").Value = "A"
$ws.Range("J1").Value = "A"
$ws.Range("K1").Value = 10

# Row 2
$ws.Range("B2").Value = 10
$ws.Range("F2").Value = "B"
$ws.Range("G2").Value = "B"
$ws.Range("H2").Value = "B"
$ws.Range("I2").Value = "B"
$ws.Range("J2").Value = "B"
$ws.Range("K2").Value = 20

# Row 3
$ws.Range("B3").Value = 10
$ws.Range("F3").Value = "C"
$ws.Range("G3").Value = "C"
$ws.Range("H3").Value = "C"
$ws.Range("I3").Value = "C"
$ws.Range("J3").Value = "C"
$ws.Range("K3").Value = 20

# Row 4
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "E"
$ws.Range("D4").Value = "E"
$ws.Range("E4").Value = "E"
$ws.Range("F4").Value = "E"
$ws.Range("G4").Value = "E"
$ws.Range("H4").Value = "E"
$ws.Range("I4").Value = "E"
$ws.Range("J4").Value = "E"
$ws.Range("K4").Value = 20
